$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Cell value updates (fill some previously-missing values, clear some others) ---

# Row 3 (RM 8): C3 filled in
$ws.Range("C3").Value = 11.2

# Row 4 (RM 9): D4 cleared
$ws.Range("D4").ClearContents()

# Row 5 (RM 14): C5 cleared
$ws.Range("C5").ClearContents()

# Row 6 (RM 21): F6 filled in
$ws.Range("F6").Value = 16.43

# Row 9 (RM 42): D9 filled in
$ws.Range("D9").Value = -14.5

# Row 10 (RM 52 a): D10 filled in
$ws.Range("D10").Value = -14.7

# Row 12 (RM 81): F12 cleared
$ws.Range("F12").ClearContents()

# Row 14 (RM 90): F14 filled in
$ws.Range("F14").Value = 17.76

# Row 17 (RM 116): D17 cleared, F17 filled in
$ws.Range("D17").ClearContents()
$ws.Range("F17").Value = 17.78

# Row 18 (RM 120): D18 cleared
$ws.Range("D18").ClearContents()

# Row 19 (RM 125): F19 filled in
$ws.Range("F19").Value = 17.81

# Row 20 (RM 134): F20 cleared
$ws.Range("F20").ClearContents()

# Row 21 (RM 135): C21 filled in
$ws.Range("C21").Value = 12.7

# Row 23 (RM 140): C23 cleared, F23 cleared
$ws.Range("C23").ClearContents()
$ws.Range("F23").ClearContents()

# Row 25 (RM 145): F25 cleared
$ws.Range("F25").ClearContents()

# Row 34 (SC 193): C34 filled in
$ws.Range("C34").Value = 10.5

# Row 29 (SC 101): F29 filled in
$ws.Range("F29").Value = 17

# Row 30 (SC 105): F30 filled in
$ws.Range("F30").Value = 17.44

# --- Row removals ---
# Remove the "RM 232" row (originally row 26) and the "SC 92" row
# (originally row 28, which becomes row 27 once row 26 is removed).
$ws.Rows("26").Delete()
$ws.Rows("27").Delete()
